# Tijdschriftformulier.xlsx - Logboek bijwerken + nieuwe sharedStrings voor
# "toevoegen aan winkelwagen" probleem.

$wb = $excel.ActiveWorkbook

$activiteit  = "Probleem oplossen toevoegen aan winkelwagen"
$opmerking   = "Probleem: Het toevoegen werkt niet. Dit is opgelost. Daarna was er een nieuw probleem: Je kan niet meer dan 1 item toevoegen. Dit kwam omdat het aantal een string was en we konden er geen integer van maken. Dit is bijna opgelost, alleen uit het functioneel ontwerp is gebleken dat we dit alleen moeten doen bij het aanpassen van de winkelmand."

# Sheets that need the new logboek entry on row 45 (last row of the
# underlying table, A9:D45): "P4 - Ivar" en "P6 - Jasper"
$sheetNames = @("P4 - Ivar", "P6 - Jasper")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A45").Value = $activiteit
    $ws.Range("B45").Value = 44146
    $ws.Range("B45").NumberFormat = "mm-dd-yy"
    $ws.Range("C45").Value = 120
    $ws.Range("D45").Value = $opmerking
}

# Re-select cells to match the saved selection state in the workbook.
$wsIvar = $wb.Worksheets.Item("P4 - Ivar")
$wsIvar.Activate()
$wsIvar.Range("A48").Select()

$wsJasper = $wb.Worksheets.Item("P6 - Jasper")
$wsJasper.Activate()
$wsJasper.Range("A45:D45").Select()
$wsJasper.Range("D45").Activate()

$wb.Save()
